$d = $word.ActiveDocument

# 1) Merge "... è una relazione " + "OneToMany" + " dato che ... la relazione " into one run
$d.Content.Find.Execute(
    " è una relazione OneToMany dato che un utente può avere molti prestiti, di conseguenza la relazione ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " è una relazione OneToMany dato che un utente può avere molti prestiti, di conseguenza la relazione ",
    2) | Out-Null

# 2) Merge " sarà " + "ManyToOne" + "." into one run
$d.Content.Find.Execute(
    " sarà ManyToOne.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " sarà ManyToOne.",
    2) | Out-Null

# 3) Merge "La relazione Prestito-Catalogo è di tipo " + "OneToMany" + " dato che ... prestiti." into one run
$d.Content.Find.Execute(
    "La relazione Prestito-Catalogo è di tipo OneToMany dato che un elemento del catalogo può essere associato a molti prestiti.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "La relazione Prestito-Catalogo è di tipo OneToMany dato che un elemento del catalogo può essere associato a molti prestiti.",
    2) | Out-Null

# 4) Append a new paragraph with the N.B. note at the end of the document
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "N.B. Come ci hai insegnato consiglio di lanciare il MainCreate con il persistence su create per avere una struttura iniziale come la mia, quindi i db popolati con libri riviste 4 utenti e molti prestiti, poi si passa al update in persitence e quindi si lancia il MainMenu che per errore ho lasciato nel package entity, quindi girare il menu utente a piacere."
